$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "simplesentiment" label to "simpleconnotation" (row 10, column A).
# Excel's shared-string table drops the now-unused "simplesentiment" entry and
# appends "simpleconnotation" as a new shared string; every other cell that
# referenced a shared string after the old slot shifts down automatically.
$ws.Range("A10").Value = "simpleconnotation"

# Update the view: scroll so row 4 becomes the top visible row (column A stays
# the left-most visible column), then select E10 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E10").Select()
